# Add one more QUERY test case (row 13) to Sheet1, mirroring the existing
# "SkillsForLevelOfEnglish" tests in rows 10-11 but passing the numeric
# arguments as a cell range (C13:E13) instead of as literal numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Supporting values for the range-based QUERY call.
$ws.Range("C13").Value = 3.8
$ws.Range("D13").Value = 2.4
$ws.Range("E13").Value = 4.4000000000000004

# Formula cell: INDEX(QUERY("SkillsForLevelOfEnglish","A11",C13:E13),2,1)
$ws.Range("A13").Formula = '=INDEX(QUERY("SkillsForLevelOfEnglish","A11",C13:E13),2,1)'

# Expected-name column, same shared value used in rows 7-11 ("Lucasz").
$ws.Range("B13").Value = "Lucasz"

# Selection follows the newly entered row, as in the target workbook.
$ws.Range("B13").Select()
